$wb = $excel.ActiveWorkbook

# Sheet "展览" and "全部类型" both contain the same event listing data.
# Row 2 (丽水·CCAC动漫游戏嘉年华) "想去人数" (F2) increments from 228 to 229.
# Row 4 (丽水·第三届HP国风动漫游戏嘉年华) "想去人数" (F4) increments from 151 to 152.

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 229
    $ws.Range("F4").Value = 152
}
